$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.986.32'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '2.294.40'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''300.42'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '''99.27'
$ws.Range("E6").Value = '  +2.29%  '
$ws.Range("E7").Value = '  -1.10%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D10").Value = '''36.14'
$ws.Range("E10").Value = '  +8.09%  '
$ws.Range("D11").Value = '''0.0789'
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("D12").Value = '''0.117'
$ws.Range("E12").Value = '  +1.07%  '
$ws.Range("D13").Value = '''18.06'
$ws.Range("E13").Value = '  +7.63%  '
$ws.Range("E14").Value = '  +2.24%  '
$ws.Range("D15").Value = '2.651.56'
$ws.Range("D16").Value = '2.294.94'
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").Value = '42.879.33'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").Value = '''12.53'
$ws.Range("E19").Value = '  +8.59%  '
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("D22").Value = '''67.74'
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("D23").Value = '''235.85'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +9.85%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '''2.43'
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").Value = '''24.93'
$ws.Range("E27").Value = '  +2.57%  '
$ws.Range("D28").Value = '''2.35'
$ws.Range("E28").Value = '  +14.99%  '
$ws.Range("D29").Value = '''34.49'
$ws.Range("E29").Value = '  +2.26%  '
$ws.Range("D30").Value = '''167.20'
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = '''5.00'
$ws.Range("E33").Value = '  +1.60%  '
$ws.Range("D34").Value = '''17.63'
$ws.Range("E34").Value = '  +3.86%  '
$ws.Range("D35").Value = '''4.61'
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("D37").Value = '''0.0688'
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("E38").Value = '  -0.51%  '
$ws.Range("E39").Value = '  +1.88%  '
$ws.Range("E40").Value = '  -0.49%  '
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  -1.41%  '
$ws.Range("D43").Value = '''0.0292'
$ws.Range("E43").Value = '  +4.04%  '
$ws.Range("D44").Value = '1.971.39'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").Value = '''10.18'
$ws.Range("E45").Value = '  +3.40%  '
$ws.Range("E46").Value = '  +1.79%  '
$ws.Range("D47").Value = '''17.44'
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("E48").Value = '  +4.35%  '
$ws.Range("E49").Value = '  +3.84%  '
$ws.Range("D50").Value = '2.518.79'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").Value = '''70.69'
$ws.Range("E51").Value = '  +1.06%  '
